# Planner update: refresh progress status for several Data Structure /
# Concept rows, and record that flashcards have been partially created
# for Bit Manipulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Status column (A) updates -------------------------------------------
# Trees, Graphs: previously N/a -> Done
$ws.Range("A5").Value = "Done"
$ws.Range("A7").Value = "Done"

# Tries: previously N/a -> In Prog
$ws.Range("A6").Value = "In Prog"

# Stacks, Queues: previously In Prog -> Done
$ws.Range("A8").Value = "Done"
$ws.Range("A9").Value = "Done"

# Big O (Concept): previously N/a -> Review
$ws.Range("A21").Value = "Review"

# Bit Manipulation (Concept): previously N/a -> In Prog, as flashcard
# work has now partially started on it.
$ws.Range("A22").Value = "In Prog"
$ws.Range("E22").Value = "In Prog"

# Leave the cursor on the cell that was last edited, matching the
# workbook's saved view state.
$ws.Range("E22").Select()
